$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.254.74"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.324.92"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.03"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.54"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("E8").Value = "  -1.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.321.48"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.336"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.59"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.736.33"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.230.83"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.320.13"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.60"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.09"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.15"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.10"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.37"
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.27"
$ws.Range("E29").Value = "  +10.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.21"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0730"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.05"
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.36"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.00"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "318.78"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.04"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.52"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.13"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.51"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0945"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.17"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0495"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0216"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.94"
$ws.Range("E51").Value = "  -0.70%  "
